$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force updated Price/Volume cells to remain Text (matches source data,
# which stores these as inline strings, e.g. "27.451.70", "0.00001030", "118.10")
# so Excel does not auto-coerce numeric-looking values into numbers.

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.451.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.92%  "

# Row 4: TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.50%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6: USDC
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4578"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.05%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3828"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.76%  "

# Row 9: OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.45"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.39%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07919"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.89%  "

# Row 11: Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9682"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.22%  "

# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.06%  "

# Row 13: WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.849.98"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.95%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.872"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.87%  "

# Row 15: Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.048"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.53%  "

# Row 16: BinanceUSD
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.60%  "

# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.16%  "

# Row 18: TRON
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06664"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001030"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.14%  "

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.87%  "

# Row 21: Dai
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.62%  "

# Row 22: WrappedBTC
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.442.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.98%  "

# Row 23: Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.336"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.28%  "

# Row 24: Cosmos
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.95%  "

# Row 25: Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26: WrappedliquidstakedEther2.0
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.061.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.28%  "

# Row 27: Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.28"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.50%  "

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.36"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29: LidoDAOToken
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.98%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.78%  "

# Row 31: BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.10"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.31%  "

# Row 32: ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9462"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.57%  "

# Row 33: Stellar
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09287"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.58%  "

# Row 34: HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.579"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.70%  "

# Row 35: Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.235"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36: ARBITRUM
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.48%  "

# Row 37: Hedera
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05934"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.66%  "

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02195"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.77%  "

# Row 39: TrustWalletToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.158"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.51%  "

# Row 40: FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.996"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.04%  "

# Row 41: TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5782"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.20%  "

# Row 42: Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1835"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.40%  "

# Row 43: Aptos
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.77%  "

# Row 44: WEMIXTOKEN
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.278"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.54%  "

# Row 45: Decentraland
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5481"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.52%  "

# Row 46: EnergySwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.01"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.25%  "

# Row 47: NEARProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.865"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.68%  "

# Row 48: Cronos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06644"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.82%  "

# Row 49: Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.80"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.45%  "

# Row 50: EOS
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.01%  "

# Row 51: PaxDollar
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.54%  "

